# Insert a new "name" column (C) into the worksheet, pushing the existing
# "geometry" column from C to D, and populate it with station names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (geometry shifts to D).
$ws.Columns.Item(3).Insert()

# Header
$ws.Cells.Item(1, 3).Value = "name"

# Station names, in row order (rows 2-47).
$names = @(
  "Diemen Zuid",
  "Amstelveenseweg",
  "Jan van Galenstraat",
  "Amsterdam Holendrecht",
  "Van der Madeweg",
  "Overamstel",
  "Halfweg-Zwanenburg",
  "De Vlugtlaan",
  "Amsterdam Sloterdijk",
  "Amsterdam Bijlmer ArenA",
  "Amsterdam Amstel",
  "Gein",
  "Reigersbos",
  "Gaasperplas",
  "Ganzenhoef",
  "Kraaiennest",
  "Venserpolder",
  "Duivendrecht",
  "Amsterdam RAI",
  "Strandvliet",
  "Bullewijk",
  "Heemstedestraat",
  "Verrijn Stuartweg",
  "Spaklerweg",
  "Wibautstraat",
  "Nieuwmarkt",
  "Henk Sneevlietweg",
  "Amsterdam Lelylaan",
  "Isolatorweg",
  "Postjesweg",
  "Zaandam",
  "Amsterdam Centraal",
  "Amsterdam Muiderpoort",
  "Amsterdam Zuid",
  "Zuid",
  "Station Sloterdijk",
  "Station RAI",
  "Centraal Station",
  "Europaplein",
  "De Pijp",
  "Vijzelgracht",
  "Noorderpark",
  "Noord",
  "Amsterdam Science Park",
  "Diemen",
  "Centraal Station"
)

$row = 2
foreach ($name in $names) {
  $ws.Cells.Item($row, 3).Value = $name
  $row = $row + 1
}
